# rajout de tests matt
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlCenter = -4108

# B column (name) then C column (description) are entered in the same order
# the author typed them -- row by row, except the "Bases" pair (rows 66/67)
# where both names were typed before either description.

$ws.Cells.Item(62, 2).Value = "Plateau"
$ws.Cells.Item(62, 2).HorizontalAlignment = $xlCenter

$ws.Cells.Item(63, 2).Value = "testDimensionsDebut"
$ws.Cells.Item(63, 2).HorizontalAlignment = $xlCenter
$ws.Cells.Item(63, 3).Value = "validé si le plateau est bien de la dimension entrée en paramètres au début de la partie"

$ws.Cells.Item(64, 2).Value = "testDimensionsCours"
$ws.Cells.Item(64, 2).HorizontalAlignment = $xlCenter
$ws.Cells.Item(64, 3).Value = "validé si le plateau fait toujours la même dimension au cours de cette même partie"

$ws.Cells.Item(65, 2).Value = "testObstacles"
$ws.Cells.Item(65, 2).HorizontalAlignment = $xlCenter
$ws.Cells.Item(65, 3).Value = "validé si le plateau contient bien le bon pourcentage d'obstacles entré en paramètre au début de la partie"

$ws.Cells.Item(66, 2).Value = "testBasesDebut"
$ws.Cells.Item(66, 2).HorizontalAlignment = $xlCenter
$ws.Cells.Item(67, 2).Value = "testBasesCours"
$ws.Cells.Item(67, 2).HorizontalAlignment = $xlCenter
$ws.Cells.Item(66, 3).Value = "validé si les bases sont situées aux extrémités d'une diagonale au début de la partie"
$ws.Cells.Item(67, 3).Value = "validé si les bases sont situées sur les même cases qu'au début de cette même partie"

$ws.Cells.Item(68, 2).Value = "testSortieImpossibleDebut"
$ws.Cells.Item(68, 2).HorizontalAlignment = $xlCenter
$ws.Cells.Item(68, 3).Value = "validé si la sortie de plateau par chaque côté est impossible au début de la partie"

$ws.Cells.Item(69, 2).Value = "testSortieImpossibleTireur"
$ws.Cells.Item(69, 2).HorizontalAlignment = $xlCenter
$ws.Cells.Item(69, 3).Value = "validé si un tireur ne peut pas sortir du plateau en cours de partie"

$ws.Cells.Item(70, 2).Value = "testSortieImpossiblePiegeur"
$ws.Cells.Item(70, 2).HorizontalAlignment = $xlCenter
$ws.Cells.Item(70, 3).Value = "validé si un piégeur ne peut pas sortir du plateau en cours de partie"

$ws.Cells.Item(71, 2).Value = "testSortieImpossibleChar"
$ws.Cells.Item(71, 2).HorizontalAlignment = $xlCenter
$ws.Cells.Item(71, 3).Value = "validé si un char ne peut pas sortir du plateau en cours de partie"

$ws.Application.ActiveWindow.ScrollRow = 48
$ws.Range("B72").Select()
